$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(4, 10).Value = 1.05
$ws.Cells.Item(4, 11).Value = 11
$ws.Cells.Item(4, 12).Value = 1.25
$ws.Cells.Item(4, 13).Value = 3.75
$ws.Cells.Item(4, 21).Value = 9
$ws.Cells.Item(4, 26).Value = 11
$ws.Cells.Item(4, 30).Value = 201
$ws.Cells.Item(4, 33).Value = 15
$ws.Cells.Item(12, 7).Value = 3.5
$ws.Cells.Item(12, 8).Value = 3.7
$ws.Cells.Item(12, 9).Value = 1.95
$ws.Cells.Item(12, 21).Value = 17
$ws.Cells.Item(12, 22).Value = 12
$ws.Cells.Item(12, 24).Value = 29
$ws.Cells.Item(12, 29).Value = 51
$ws.Cells.Item(12, 32).Value = 9
$ws.Cells.Item(13, 11).Value = 13
$ws.Cells.Item(13, 12).Value = 1.25
$ws.Cells.Item(13, 13).Value = 3.75
$ws.Cells.Item(13, 14).Value = 1.8
$ws.Cells.Item(13, 15).Value = 2
$ws.Cells.Item(14, 7).Value = 2.2
$ws.Cells.Item(14, 9).Value = 3
$ws.Cells.Item(14, 14).Value = 1.73
$ws.Cells.Item(14, 15).Value = 2.08
$ws.Cells.Item(19, 7).Value = 3.4
$ws.Cells.Item(19, 9).Value = 2.2
$ws.Cells.Item(19, 12).Value = 1.35
$ws.Cells.Item(19, 13).Value = 2.7
$ws.Cells.Item(19, 14).Value = 2.02
$ws.Cells.Item(19, 15).Value = 1.62
$ws.Cells.Item(19, 16).Value = 1.47
$ws.Cells.Item(19, 17).Value = 2.32
$ws.Cells.Item(19, 18).Value = 1.75
$ws.Cells.Item(19, 19).Value = 1.85
$ws.Cells.Item(19, 20).Value = 8.75
$ws.Cells.Item(19, 26).Value = 7.8
$ws.Cells.Item(19, 28).Value = 13.5
$ws.Cells.Item(19, 31).Value = 7.2
$ws.Cells.Item(19, 32).Value = 10.75
$ws.Cells.Item(19, 34).Value = 23
$ws.Cells.Item(19, 35).Value = 18
$ws.Cells.Item(19, 36).Value = 29
$ws.Cells.Item(20, 7).Value = 2.35
$ws.Cells.Item(20, 9).Value = 2.95
$ws.Cells.Item(20, 12).Value = 1.39
$ws.Cells.Item(20, 20).Value = 6.8
$ws.Cells.Item(20, 21).Value = 10.75
$ws.Cells.Item(20, 22).Value = 9.5
$ws.Cells.Item(20, 23).Value = 24
$ws.Cells.Item(20, 31).Value = 7.9
$ws.Cells.Item(20, 32).Value = 14
$ws.Cells.Item(20, 33).Value = 11
$ws.Cells.Item(20, 34).Value = 37
$ws.Cells.Item(20, 35).Value = 28
$ws.Cells.Item(31, 7).Value = 5.5
$ws.Cells.Item(31, 9).Value = 1.45
$ws.Cells.Item(31, 10).Value = 26
$ws.Cells.Item(31, 11).Value = 1.02
$ws.Cells.Item(31, 25).Value = 34
$ws.Cells.Item(31, 32).Value = 11
$ws.Cells.Item(33, 7).Value = 1.57
$ws.Cells.Item(33, 8).Value = 4.5
$ws.Cells.Item(33, 9).Value = 4.5
$ws.Cells.Item(33, 10).Value = 26
$ws.Cells.Item(33, 11).Value = 1.02
$ws.Cells.Item(33, 16).Value = 1.2
$ws.Cells.Item(33, 17).Value = 4.33
$ws.Cells.Item(33, 18).Value = 1.44
$ws.Cells.Item(33, 19).Value = 2.63
$ws.Cells.Item(33, 20).Value = 13
$ws.Cells.Item(33, 21).Value = 12
$ws.Cells.Item(33, 26).Value = 26
$ws.Cells.Item(33, 27).Value = 10
$ws.Cells.Item(33, 30).Value = 81
$ws.Cells.Item(33, 31).Value = 23
$ws.Cells.Item(33, 33).Value = 17
$ws.Cells.Item(35, 9).Value = 2.1
$ws.Cells.Item(35, 10).Value = 1.02
$ws.Cells.Item(35, 11).Value = 11
$ws.Cells.Item(35, 24).Value = 23
$ws.Cells.Item(35, 25).Value = 29
$ws.Cells.Item(35, 26).Value = 11
$ws.Cells.Item(35, 33).Value = 9.5
$ws.Cells.Item(36, 7).Value = 2.9
$ws.Cells.Item(36, 8).Value = 3.75
$ws.Cells.Item(36, 9).Value = 2.15
$ws.Cells.Item(36, 10).Value = 1.03
$ws.Cells.Item(36, 11).Value = 10
$ws.Cells.Item(36, 20).Value = 12
$ws.Cells.Item(40, 7).Value = 3.1
$ws.Cells.Item(40, 8).Value = 3.7
$ws.Cells.Item(40, 9).Value = 2.15
$ws.Cells.Item(40, 16).Value = 1.33
$ws.Cells.Item(40, 17).Value = 3.25
$ws.Cells.Item(40, 18).Value = 1.62
$ws.Cells.Item(40, 19).Value = 2.2
$ws.Cells.Item(40, 20).Value = 12
$ws.Cells.Item(40, 21).Value = 17
$ws.Cells.Item(40, 22).Value = 11
$ws.Cells.Item(40, 25).Value = 26
$ws.Cells.Item(40, 26).Value = 13
$ws.Cells.Item(40, 32).Value = 12
$ws.Cells.Item(40, 34).Value = 21
$ws.Cells.Item(40, 35).Value = 17
$ws.Cells.Item(40, 36).Value = 23
$ws.Cells.Item(43, 7).Value = 2.05
$ws.Cells.Item(43, 8).Value = 3.5
$ws.Cells.Item(43, 9).Value = 3.4
$ws.Cells.Item(43, 20).Value = 8
$ws.Cells.Item(43, 21).Value = 10
$ws.Cells.Item(43, 22).Value = 9
$ws.Cells.Item(43, 23).Value = 19
$ws.Cells.Item(43, 24).Value = 17
$ws.Cells.Item(43, 27).Value = 6.5
$ws.Cells.Item(43, 29).Value = 41
$ws.Cells.Item(43, 32).Value = 17
$ws.Cells.Item(43, 33).Value = 12
$ws.Cells.Item(43, 35).Value = 26
$ws.Cells.Item(45, 12).Value = 1.18
$ws.Cells.Item(45, 13).Value = 4.5
$ws.Cells.Item(45, 14).Value = 1.62
$ws.Cells.Item(45, 15).Value = 2.25
$ws.Cells.Item(46, 14).Value = 1.2
$ws.Cells.Item(46, 15).Value = 4.33
